$d = $word.ActiveDocument
$p2 = $d.Paragraphs(2)
$r = $p2.Range
$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>P</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>rivacy is something that is becoming  more important to all of us.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>A long time ago,it wasn</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>’</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>t so important.Our world was small.We knew and trust</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>ed</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> everyone around us.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>We rarely had to give people our private information. Today, with the Internet and modern communications, privacy is a hot topic. It is now possible to steal someone’s identity. There are even websites that sell people’s private details. We have to be really careful with what we put online. Especially on social network sites. The private info you put on these is there for the whole world to see. It’s easy to search for someone’s name in Google and find out all kinds of stuff. I think people need to think more about privacy and be more careful.</w:t>
      </w:r>
    </w:p>
'@
$r.InsertXML($xml)

$style = $d.Styles.Add("apple-converted-space", 2)
$style.BaseStyle = $d.Styles("a0")
